$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (row 1): B1:E1
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Update row 2 values (B2:E2)
$ws.Range("B2").Value = 277.55526040324696
$ws.Range("C2").Value = 258.27178096320824
$ws.Range("D2").Value = 278.09091106517894
$ws.Range("E2").Value = 257.47125743493308

# Update row 3 values (B3:E3)
$ws.Range("B3").Value = 304.71671936823464
$ws.Range("C3").Value = 259.59248883217685
$ws.Range("D3").Value = 310.08155626722527
$ws.Range("E3").Value = 254.3701899394438

# Update the selected range to match new selection
$ws.Range("B1:E3").Select()
